$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ebullition")

$ws.Range("B4").Value = "07152016.BWL.CAL1.0.1.UNK"
$ws.Range("B24").Value = "09022016CAL1"
$ws.Range("B19").Value = "08122016.CAL1"
$ws.Range("B27").Value = "09152016.DIL.S06.13.5'.UNK"

$ws.Range("B28").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
